$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Sound"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "Death sound"
$ws.Range("D9").Value = "Plays when player runs out of health"

$ws.Range("D9").Select()
